$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

$ws.Range("A4").Value = 264163752684.3103
$ws.Range("B4").Value = 39223134532.18238
$ws.Range("C4").Value = 3542495631.375376
$ws.Range("D4").Value = 215902978315.8318
$ws.Range("E4").Value = 5495144204.920725
$ws.Range("F4").Value = 255126112848.0142
$ws.Range("G4").Value = 9037639836.296101
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 264163752684.3103
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 29473.37700009346
$ws.Range("W4").Value = 263292262834.5643
$ws.Range("X4").Value = 264163752684.3103
$ws.Range("Y4").Value = 871489849.7460327
$ws.Range("Z4").Value = "optimal"
$ws.Range("AA4").Value = "costs_emissionlimit"
$ws.Range("AB4").Value = -1
$ws.Range("AC4").Value = -1
$ws.Range("AD4").Value = 1
$ws.Range("AE4").Value = "GreenFieldHydro_Island"
$ws.Range("AF4").Value = "rawResults\20251119115745_GreenFieldHydro_Island-1"
